$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove the "Profiling" marker (column E) from rows 2, 3 and 4 ---
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()

# --- 2. Add an empty, but styled (same fill as other column C cells), cell in C5 and C6 ---
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C6").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()

# --- 3. Insert a new row at position 15 for "ActualizarPruebasMedicas" ---
# (this shifts old row 15 "RegistrarPetType" -> 16 and old row 16 "ActualizarPetType" -> 17,
#  carrying their formatting down automatically)
$ws.Rows("15:15").Insert()

$ws.Range("A15").Value = "ActualizarPruebasMedicas"
$ws.Range("B15").Value = 500
$ws.Range("C15").Value = 3000
$ws.Range("E15").Value = "Profiling"

# --- 4. Append four new rows (18-21) describing the new "Tratamiento" user stories ---
# Row 18 and 19 follow the same pattern as row 14 (A, B, C, E columns with the "Profiling" tag)
$ws.Range("A14").Copy()
$ws.Range("A18").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A19").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$ws.Range("B19").PasteSpecial(-4122)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E19").PasteSpecial(-4122)

# Row 20 follows the plain A, B, C pattern (no E tag), like row 16
$ws.Range("A16").Copy()
$ws.Range("A20").PasteSpecial(-4122)
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("B16").Copy()
$ws.Range("B20").PasteSpecial(-4122)
$ws.Range("C16").Copy()
$ws.Range("C20").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("A18").Value = "CrearTratamiento"
$ws.Range("B18").Value = 100
$ws.Range("C18").Value = 1500
$ws.Range("E18").Value = "Profiling"

$ws.Range("A19").Value = "ActualizarTratamiento"
$ws.Range("B19").Value = 100
$ws.Range("C19").Value = 1500
$ws.Range("E19").Value = "Profiling"

$ws.Range("A20").Value = "EliminarHitorialTratamiento"
$ws.Range("B20").Value = 200
$ws.Range("C20").Value = 1500

$ws.Range("A21").Value = "ListarTratamientos"

# --- 5. Update the saved selection to match the author's final cursor position ---
[void]$ws.Range("B21").Select()
